$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue {
    param($cellRange, [string]$text)
    $cellRange.NumberFormat = '@'
    $cellRange.Value = $text
    $cellRange.Style = 'Normal'
}

Set-TextValue $ws.Range("D2") "29.538.06"
Set-TextValue $ws.Range("E2") "  +0.16%  "
Set-TextValue $ws.Range("D3") "1.911.52"
Set-TextValue $ws.Range("E3") "  -0.10%  "
Set-TextValue $ws.Range("D4") "1.005"
Set-TextValue $ws.Range("E4") "  +0.54%  "
Set-TextValue $ws.Range("D5") "326.40"
Set-TextValue $ws.Range("E5") "  -0.40%  "
Set-TextValue $ws.Range("D6") "1.005"
Set-TextValue $ws.Range("E6") "  +0.48%  "
Set-TextValue $ws.Range("D7") "0.4845"
Set-TextValue $ws.Range("E7") "  +2.09%  "
Set-TextValue $ws.Range("E8") "  -0.42%  "
Set-TextValue $ws.Range("E9") "  +1.54%  "
Set-TextValue $ws.Range("D10") "1.014"
Set-TextValue $ws.Range("E10") "  +0.34%  "
Set-TextValue $ws.Range("D11") "23.52"
Set-TextValue $ws.Range("E11") "  +4.95%  "
Set-TextValue $ws.Range("D12") "1.902.41"
Set-TextValue $ws.Range("E12") "  -0.44%  "
Set-TextValue $ws.Range("D13") "6.031"
Set-TextValue $ws.Range("E13") "  +1.28%  "
Set-TextValue $ws.Range("D14") "7.112"
Set-TextValue $ws.Range("E14") "  -0.73%  "
Set-TextValue $ws.Range("D15") "90.42"
Set-TextValue $ws.Range("E15") "  +1.17%  "
Set-TextValue $ws.Range("B16") "BinanceUSD"
Set-TextValue $ws.Range("C16") "https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd"
Set-TextValue $ws.Range("D16") "1.006"
Set-TextValue $ws.Range("E16") "  +0.63%  "
Set-TextValue $ws.Range("B17") "TRON"
Set-TextValue $ws.Range("C17") "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
Set-TextValue $ws.Range("D17") "0.06763"
Set-TextValue $ws.Range("E17") "  +2.52%  "
Set-TextValue $ws.Range("D18") "0.00001042"
Set-TextValue $ws.Range("E18") "  +1.31%  "
Set-TextValue $ws.Range("E19") "  -0.32%  "
Set-TextValue $ws.Range("E20") "  +0.48%  "
Set-TextValue $ws.Range("D21") "29.549.20"
Set-TextValue $ws.Range("E21") "  +0.15%  "
Set-TextValue $ws.Range("D22") "5.614"
Set-TextValue $ws.Range("E22") "  +1.11%  "
Set-TextValue $ws.Range("D23") "11.82"
Set-TextValue $ws.Range("E23") "  +2.41%  "
Set-TextValue $ws.Range("D24") "2.166"
Set-TextValue $ws.Range("E24") "  -1.80%  "
Set-TextValue $ws.Range("D25") "2.138.81"
Set-TextValue $ws.Range("E25") "  -0.34%  "
Set-TextValue $ws.Range("D26") "154.81"
Set-TextValue $ws.Range("E26") "  +0.94%  "
Set-TextValue $ws.Range("D27") "20.07"
Set-TextValue $ws.Range("E27") "  +1.36%  "
Set-TextValue $ws.Range("D28") "6.292"
Set-TextValue $ws.Range("E28") "  +9.57%  "
Set-TextValue $ws.Range("D29") "2.109"
Set-TextValue $ws.Range("E29") "  -1.30%  "
Set-TextValue $ws.Range("D30") "119.82"
Set-TextValue $ws.Range("E30") "  +1.98%  "
Set-TextValue $ws.Range("D31") "1.034"
Set-TextValue $ws.Range("E31") "  -3.17%  "
Set-TextValue $ws.Range("D32") "0.09563"
Set-TextValue $ws.Range("E32") "  +0.10%  "
Set-TextValue $ws.Range("D33") "5.539"
Set-TextValue $ws.Range("E33") "  +2.71%  "
Set-TextValue $ws.Range("D34") "1.398"
Set-TextValue $ws.Range("E34") "  -1.76%  "
Set-TextValue $ws.Range("E35") "  -0.61%  "
Set-TextValue $ws.Range("D36") "0.02270"
Set-TextValue $ws.Range("E36") "  +0.64%  "
Set-TextValue $ws.Range("D37") "0.06121"
Set-TextValue $ws.Range("E37") "  +0.41%  "
Set-TextValue $ws.Range("D38") "1.174"
Set-TextValue $ws.Range("E38") "  +0.16%  "
Set-TextValue $ws.Range("D39") "0.5957"
Set-TextValue $ws.Range("E39") "  +1.09%  "
Set-TextValue $ws.Range("D40") "7.945"
Set-TextValue $ws.Range("D41") "10.68"
Set-TextValue $ws.Range("E41") "  +5.33%  "
Set-TextValue $ws.Range("D42") "0.1857"
Set-TextValue $ws.Range("E42") "  +0.88%  "
Set-TextValue $ws.Range("D43") "2.444"
Set-TextValue $ws.Range("E43") "  +1.35%  "
Set-TextValue $ws.Range("D44") "1.284"
Set-TextValue $ws.Range("E44") "  -1.34%  "
Set-TextValue $ws.Range("D45") "0.07715"
Set-TextValue $ws.Range("E45") "  -1.01%  "
Set-TextValue $ws.Range("D46") "12.41"
Set-TextValue $ws.Range("E46") "  +1.60%  "
Set-TextValue $ws.Range("D47") "0.5579"
Set-TextValue $ws.Range("E47") "  +0.51%  "
Set-TextValue $ws.Range("D48") "1.957"
Set-TextValue $ws.Range("E48") "  +1.22%  "
Set-TextValue $ws.Range("D49") "115.01"
Set-TextValue $ws.Range("E49") "  +1.30%  "
Set-TextValue $ws.Range("D50") "72.77"
Set-TextValue $ws.Range("E50") "  +1.72%  "
Set-TextValue $ws.Range("D51") "1.054"
Set-TextValue $ws.Range("E51") "  +2.13%  "
